$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 9).Value = "sd"
$ws.Cells.Item(2, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(4, 9).Value = "sv"
$ws.Cells.Item(4, 10).Value = "Statement-opinion"
$ws.Cells.Item(11, 9).Value = "sv"
$ws.Cells.Item(11, 10).Value = "Statement-opinion"
$ws.Cells.Item(27, 9).Value = "aa"
$ws.Cells.Item(27, 10).Value = "Agree/Accept"
$ws.Cells.Item(32, 9).Value = "b"
$ws.Cells.Item(32, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(37, 9).Value = "sd"
$ws.Cells.Item(37, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(45, 9).Value = "aa"
$ws.Cells.Item(45, 10).Value = "Agree/Accept"
$ws.Cells.Item(53, 9).Value = "aa"
$ws.Cells.Item(53, 10).Value = "Agree/Accept"
$ws.Cells.Item(55, 9).Value = "aa"
$ws.Cells.Item(55, 10).Value = "Agree/Accept"
$ws.Cells.Item(63, 9).Value = "aa"
$ws.Cells.Item(63, 10).Value = "Agree/Accept"
$ws.Cells.Item(64, 9).Value = "ba"
$ws.Cells.Item(64, 10).Value = "Appreciation"
$ws.Cells.Item(67, 9).Value = "sd"
$ws.Cells.Item(67, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(68, 9).Value = "sd"
$ws.Cells.Item(68, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(93, 9).Value = "aa"
$ws.Cells.Item(93, 10).Value = "Agree/Accept"
$ws.Cells.Item(95, 9).Value = "sd"
$ws.Cells.Item(95, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(106, 9).Value = "ba"
$ws.Cells.Item(106, 10).Value = "Appreciation"
$ws.Cells.Item(112, 9).Value = "ba"
$ws.Cells.Item(112, 10).Value = "Appreciation"
$ws.Cells.Item(120, 9).Value = "sd"
$ws.Cells.Item(120, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(140, 9).Value = "b"
$ws.Cells.Item(140, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(144, 9).Value = "sd"
$ws.Cells.Item(144, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(148, 9).Value = "aa"
$ws.Cells.Item(148, 10).Value = "Agree/Accept"
$ws.Cells.Item(150, 9).Value = "sd"
$ws.Cells.Item(150, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(157, 9).Value = "aa"
$ws.Cells.Item(157, 10).Value = "Agree/Accept"
$ws.Cells.Item(158, 9).Value = "aa"
$ws.Cells.Item(158, 10).Value = "Agree/Accept"
$ws.Cells.Item(165, 9).Value = "sd"
$ws.Cells.Item(165, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(168, 9).Value = "aa"
$ws.Cells.Item(168, 10).Value = "Agree/Accept"
$ws.Cells.Item(169, 9).Value = "aa"
$ws.Cells.Item(169, 10).Value = "Agree/Accept"
$ws.Cells.Item(196, 9).Value = "aa"
$ws.Cells.Item(196, 10).Value = "Agree/Accept"
$ws.Cells.Item(197, 9).Value = "aa"
$ws.Cells.Item(197, 10).Value = "Agree/Accept"
$ws.Cells.Item(199, 9).Value = "sd"
$ws.Cells.Item(199, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(202, 9).Value = "sd"
$ws.Cells.Item(202, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(205, 9).Value = "sd"
$ws.Cells.Item(205, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(216, 9).Value = "aa"
$ws.Cells.Item(216, 10).Value = "Agree/Accept"
$ws.Cells.Item(217, 9).Value = "aa"
$ws.Cells.Item(217, 10).Value = "Agree/Accept"
$ws.Cells.Item(218, 9).Value = "ba"
$ws.Cells.Item(218, 10).Value = "Appreciation"
$ws.Cells.Item(226, 9).Value = "sv"
$ws.Cells.Item(226, 10).Value = "Statement-opinion"
$ws.Cells.Item(231, 9).Value = "sd"
$ws.Cells.Item(231, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(233, 9).Value = "%"
$ws.Cells.Item(233, 10).Value = "Uninterpretable"
$ws.Cells.Item(241, 9).Value = "sd"
$ws.Cells.Item(241, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(244, 9).Value = "sd"
$ws.Cells.Item(244, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(254, 9).Value = "aa"
$ws.Cells.Item(254, 10).Value = "Agree/Accept"
$ws.Cells.Item(256, 9).Value = "aa"
$ws.Cells.Item(256, 10).Value = "Agree/Accept"
$ws.Cells.Item(265, 9).Value = "sv"
$ws.Cells.Item(265, 10).Value = "Statement-opinion"
$ws.Cells.Item(281, 9).Value = "sd"
$ws.Cells.Item(281, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(283, 9).Value = "sd"
$ws.Cells.Item(283, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(287, 9).Value = "aa"
$ws.Cells.Item(287, 10).Value = "Agree/Accept"
$ws.Cells.Item(288, 9).Value = "b"
$ws.Cells.Item(288, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(308, 9).Value = "sd"
$ws.Cells.Item(308, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(319, 9).Value = "sd"
$ws.Cells.Item(319, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(328, 9).Value = "sd"
$ws.Cells.Item(328, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(350, 9).Value = "ba"
$ws.Cells.Item(350, 10).Value = "Appreciation"
$ws.Cells.Item(363, 9).Value = "ba"
$ws.Cells.Item(363, 10).Value = "Appreciation"

Write-Output "Applied dialog act updates"
